# Update "Recommandations" sheet: new dimension A1:G34, refreshed stats for
# sectors (rows 2-11), and re-ranked/refreshed individual stock rows (12-34),
# including newly added SICOR CI (SICC) row and several name re-orderings.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")

$recoRows = @(
    @("BRVM - CONSOMMATION DE BASE     (**)", 0, 5, 1414.26, 281.89, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL     (**)", 0, 5, 1356.72, 275.44, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 5, 998.8, 212.11, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 5, 965.45, 198.33, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 5, 830.6, 169.85, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 5, 778.66, 159.42, "🟡 Observer", "➖ Neutre"),
    @("BRVM – COMPOSITE TOTAL RETURN     (**)", 0, 5, 759.74, 154.95, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 5, 656.79, 138.94, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 5, 615.66, 127.92, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 5, 504.88, 103.05, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 3, 0, 22.38, 7.46, "🟢 Achat", "✅ Renforcer"),
    @("SERVAIR ABIDJAN CI (ABJC)", 3, 0, 22.18, 7.37, "🟢 Achat", "✅ Renforcer"),
    @("SETAO CI (STAC)", 2, 1, 11.76, 7.4, "🟡 Observer", "👀 À surveiller"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 2, 1, 10.41, -2.94, "🟡 Observer", "👀 À surveiller"),
    @("SUCRIVOIRE (SCRC)", 2, 2, 10.17, -2.51, "🟡 Observer", "👀 À surveiller"),
    @("ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 1, 0, 7.45, 7.45, "🟡 Observer", "➖ Neutre"),
    @("EVIOSYS PACKAGING SIEM CI (SEMC)", 2, 1, 7.44, 7.33, "🟡 Observer", "👀 À surveiller"),
    @("NEI-CEDA CI (NEIC)", 1, 0, 6.92, 6.92, "🟡 Observer", "➖ Neutre"),
    @("TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 6.75, 6.75, "🟡 Observer", "➖ Neutre"),
    @("CORIS BANK INTERNATIONAL (CBIBF)", 1, 1, 5.95, -1.55, "🟡 Observer", "👀 À surveiller"),
    @("SICABLE CI (CABC)", 2, 2, 5.03, 7.39, "🟡 Observer", "👀 À surveiller"),
    @("SOLIBRA CI (SLBC)", 1, 1, 4.74, 7.49, "🟡 Observer", "👀 À surveiller"),
    @("SITAB CI (STBC)", 1, 1, 3.39, -1.9, "🟡 Observer", "👀 À surveiller"),
    @("ONATEL BF (ONTBF)", 1, 1, 3.19, -3.11, "🟡 Observer", "👀 À surveiller"),
    @("SMB CI (SMBC)", 1, 1, 0.3, 6.66, "🟡 Observer", "👀 À surveiller"),
    @("BANK OF AFRICA CI (BOAC)", 0, 1, -1.32, -1.32, "🟡 Observer", "➖ Neutre"),
    @("UNIWAX CI (UNXC)", 0, 1, -2.55, -2.55, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 0, 1, -4.1, -4.1, "🟡 Observer", "➖ Neutre"),
    @("BERNABE CI (BNBC)", 0, 3, -6.37, -1.58, "🔴 Vente", "⚠️ Risque de décrochage"),
    @("UNILEVER CI (UNLC)", 1, 2, -6.4, -6.86, "🟡 Observer", "👀 À surveiller"),
    @("ORANGE COTE D'IVOIRE (ORAC)", 0, 1, -6.93, -6.93, "🟡 Observer", "➖ Neutre"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 0, 2, -9.74, -4.41, "🟡 Observer", "➖ Neutre"),
    @("SICOR CI (SICC)", 0, 2, -14.14, -6.72, "🟡 Observer", "➖ Neutre")
)

$r = 2
foreach ($row in $recoRows) {
    $c = 1
    foreach ($val in $row) {
        $ws1.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Update "Top_YTD" sheet: refreshed "Progression YTD (%)" values (column B).
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ytdRows = @(
    @("BRVM - CONSOMMATION DE BASE     (**)", 82151.91),
    @("BRVM-PRINCIPAL     (**)", 70502.11),
    @("BRVM - INDUSTRIELS", 24036.64),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 21519.53),
    @("BRVM - SERVICES FINANCIERS", 13244.9),
    @("BRVM-PRESTIGE", 10836.2),
    @("BRVM – COMPOSITE TOTAL RETURN     (**)", 10051.1),
    @("BRVM - ENERGIE", 6521.43),
    @("BRVM - SERVICES PUBLICS", 5428.85),
    @("BRVM - TELECOMMUNICATIONS", 3178.46)
)

$r = 2
foreach ($row in $ytdRows) {
    $c = 1
    foreach ($val in $row) {
        $ws2.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
